# Roll the 90-day HTTPS/Non-HTTPS date window forward by one day:
#   - drop the oldest date row (2025-10-24)
#   - shift every remaining day up by one row
#   - append a new row for the newest date (2026-01-22) with zero counts
#
# Done on the "Chart" sheet (sheet1), which holds the Date / Non-HTTPS URLs /
# HTTPS URLs table in A1:C91 (row 1 = headers, rows 2-91 = one row per day).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 91

# Deleting row 2 physically shifts the existing rows 3..91 up to 2..90,
# carrying their original cell types/values/styles along untouched (no
# re-typing happens, so Excel's "looks like a date" auto-conversion never
# kicks in for the shifted date strings).
$ws.Rows.Item(2).Delete() | Out-Null

# Add the new trailing day in what is now the last row. Writing the date
# string straight into .Value/.Formula on a plain cell would make Excel
# auto-convert a "YYYY-MM-DD" looking value into a real date serial, which
# the source file does not use (dates are stored as plain text). Using a
# formula that evaluates to the text avoids that coercion ...
$ws.Cells.Item($lastRow, 1).Formula = '="2026-01-22"'
$ws.Cells.Item($lastRow, 2).Value = 0
$ws.Cells.Item($lastRow, 3).Value = 0

# ... then flatten the formula back down to a plain literal value so the
# cell is stored the same way as its neighbours (literal text, not a
# formula).
$ws.Cells.Item($lastRow, 1).Copy() | Out-Null
$ws.Cells.Item($lastRow, 1).PasteSpecial(-4163) | Out-Null
